$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- New data obtained for three existing requests (green "obtained" highlight + note) ---

# Row 6 (Sachsen-Anhalt line) - obtained data per email (21.11.22) -> I6
$ws.Range("I6").Value = "obtained data per email (21.11.22)"

# Row 3 (NRW line) - obtained via email 21.11.22 -> D3
$ws.Range("D3").Value = "obtained via email 21.11.22"

# Row 8 (Brandenburg line) - obtained data via email (link to cloud) 22.11.2022 -> F8
$ws.Range("F8").Value = "obtained data via email (link to cloud) 22.11.2022"

# Mark column A of these rows with the same "data obtained" highlight already used
# on rows 2, 9 and 13 (style index 5 in the workbook) by copying the format over.
$ws.Range("A2").Copy() | Out-Null
$ws.Range("A3").PasteSpecial(-4122) | Out-Null

$ws.Range("A2").Copy() | Out-Null
$ws.Range("A6").PasteSpecial(-4122) | Out-Null

$ws.Range("A2").Copy() | Out-Null
$ws.Range("A8").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = $false

# New column I needs a width similar to the other note columns.
$ws.Columns.Item(9).ColumnWidth = 25.33

# Update the view: move the selection to F9 (also resets the saved scroll
# position, dropping the old topLeftCell="A4").
$ws.Range("F9").Select() | Out-Null
